$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1099.9166
$ws.Range("I40").Value = 1111
$ws.Range("J40").Value = 1066.6666
$ws.Range("K40").Value = 1111
$ws.Range("L40").Value = 1066.6666
$ws.Range("M40").Value = -936
$ws.Range("N40").Value = -1416.6666
$ws.Range("H52").Value = 2500
$ws.Range("I52").Value = 2000
$ws.Range("J52").Value = 3000
$ws.Range("K52").Value = 6000
$ws.Range("L52").Value = 9000
$ws.Range("M52").Value = -5840
$ws.Range("N52").Value = -9320
$ws.Range("H70").Value = 22690
$ws.Range("I70").Value = 931.6667
$ws.Range("J70").Value = 29561.053
$ws.Range("K70").Value = 2795.0001
$ws.Range("L70").Value = 88683.159
$ws.Range("M70").Value = -2525.0001
$ws.Range("N70").Value = -89223.159
$ws.Range("H73").Value = 22690
$ws.Range("I73").Value = 931.6667
$ws.Range("J73").Value = 29561.053
$ws.Range("K73").Value = 2795.0001
$ws.Range("L73").Value = 88683.159
$ws.Range("M73").Value = -1859.0001
$ws.Range("N73").Value = -90555.159

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 276.57144
$ws.Range("I5").Value = 222.8
$ws.Range("J5").Value = 411
$ws.Range("K5").Value = 222.8
$ws.Range("L5").Value = 411
$ws.Range("M5").Value = -110.8
$ws.Range("N5").Value = -635
$ws.Range("H56").Value = 14800
$ws.Range("J56").Value = 14800
$ws.Range("L56").Value = 14800
$ws.Range("N56").Value = -16284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 276.57144
$ws.Range("I4").Value = 222.8
$ws.Range("J4").Value = 411
$ws.Range("K4").Value = 222.8
$ws.Range("L4").Value = 411
$ws.Range("M4").Value = -107.8
$ws.Range("N4").Value = -641
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 59933
$ws.Range("J28").Value = 59933
$ws.Range("L28").Value = 59933
$ws.Range("N28").Value = -60423
$ws.Range("H31").Value = 12049936
$ws.Range("I31").Value = 33334634
$ws.Range("J31").Value = 1993.849
$ws.Range("K31").Value = 33334634
$ws.Range("L31").Value = 1993.849
$ws.Range("M31").Value = -33334339
$ws.Range("N31").Value = -2583.849
$ws.Range("H34").Value = 12049936
$ws.Range("I34").Value = 33334634
$ws.Range("J34").Value = 1993.849
$ws.Range("K34").Value = 33334634
$ws.Range("L34").Value = 1993.849
$ws.Range("M34").Value = -33334432
$ws.Range("N34").Value = -2397.849
$ws.Range("H58").Value = 890.375
$ws.Range("I58").Value = 890.375
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 890.375
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -687.375
$ws.Range("H99").Value = 16688310
$ws.Range("I99").Value = 41686370
$ws.Range("J99").Value = 22937.834
$ws.Range("K99").Value = 41686370
$ws.Range("L99").Value = 22937.834
$ws.Range("M99").Value = -41684872
$ws.Range("N99").Value = -25933.834
$ws.Range("H126").Value = 16688310
$ws.Range("I126").Value = 41686370
$ws.Range("J126").Value = 22937.834
$ws.Range("K126").Value = 125059110
$ws.Range("L126").Value = 68813.50199999999
$ws.Range("M126").Value = -125056640
$ws.Range("N126").Value = -73753.50199999999
$ws.Range("H136").Value = 890.375
$ws.Range("I136").Value = 890.375
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2671.125
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -121.125
$ws.Range("H141").Value = 228790.33
$ws.Range("J141").Value = 256514.12
$ws.Range("L141").Value = 256514.12
$ws.Range("N141").Value = -266874.12

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3227.4893
$ws.Range("I68").Value = 650.76
$ws.Range("J68").Value = 6155.591
$ws.Range("K68").Value = 1952.28
$ws.Range("L68").Value = 18466.773
$ws.Range("M68").Value = -1141.28
$ws.Range("N68").Value = -20088.773
$ws.Range("H69").Value = 62501572
$ws.Range("I69").Value = 1117.25
$ws.Range("J69").Value = 83335060
$ws.Range("K69").Value = 3351.75
$ws.Range("L69").Value = 250005180
$ws.Range("M69").Value = -2540.75
$ws.Range("N69").Value = -250006802
$ws.Range("H71").Value = 3227.4893
$ws.Range("I71").Value = 650.76
$ws.Range("J71").Value = 6155.591
$ws.Range("K71").Value = 5856.84
$ws.Range("L71").Value = 55400.319
$ws.Range("M71").Value = -1800.84
$ws.Range("N71").Value = -63512.319
$ws.Range("H72").Value = 62501572
$ws.Range("I72").Value = 1117.25
$ws.Range("J72").Value = 83335060
$ws.Range("K72").Value = 10055.25
$ws.Range("L72").Value = 750015540
$ws.Range("M72").Value = -5999.25
$ws.Range("N72").Value = -750023652
$ws.Range("H82").Value = 31166.818
$ws.Range("I82").Value = 5000
$ws.Range("J82").Value = 33783.5
$ws.Range("K82").Value = 15000
$ws.Range("L82").Value = 101350.5
$ws.Range("M82").Value = -14594
$ws.Range("N82").Value = -102162.5
$ws.Range("H85").Value = 31166.818
$ws.Range("I85").Value = 5000
$ws.Range("J85").Value = 33783.5
$ws.Range("K85").Value = 15000
$ws.Range("L85").Value = 101350.5
$ws.Range("M85").Value = -13596
$ws.Range("N85").Value = -104158.5
$ws.Range("H107").Value = 383347.62
$ws.Range("I107").Value = 674969.7
$ws.Range("J107").Value = 1222.1724
$ws.Range("K107").Value = 2024909.1
$ws.Range("L107").Value = 3666.5172
$ws.Range("M107").Value = -2022989.1
$ws.Range("N107").Value = -7506.5172
$ws.Range("H122").Value = 1147
$ws.Range("I122").Value = 323.25
$ws.Range("J122").Value = 2245.3333
$ws.Range("K122").Value = 2909.25
$ws.Range("L122").Value = 20207.9997
$ws.Range("M122").Value = -459.25
$ws.Range("N122").Value = -25107.9997
$ws.Range("H132").Value = 1539.0571
$ws.Range("I132").Value = 1197.3334
$ws.Range("J132").Value = 1717.3478
$ws.Range("K132").Value = 10776.0006
$ws.Range("L132").Value = 15456.1302
$ws.Range("M132").Value = -8246.000599999999
$ws.Range("N132").Value = -20516.1302

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3447.842
$ws.Range("I102").Value = 3827.7144
$ws.Range("J102").Value = 2384.2
$ws.Range("K102").Value = 3827.7144
$ws.Range("L102").Value = 2384.2
$ws.Range("M102").Value = -2205.7144
$ws.Range("N102").Value = -5628.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3271480
$ws.Range("I40").Value = 4632509
$ws.Range("J40").Value = 5009.8
$ws.Range("K40").Value = 4632509
$ws.Range("L40").Value = 5009.8
$ws.Range("M40").Value = -4632373
$ws.Range("N40").Value = -5281.8
$ws.Range("H82").Value = 4417.278
$ws.Range("I82").Value = 1457.4546
$ws.Range("J82").Value = 9068.429
$ws.Range("K82").Value = 1457.4546
$ws.Range("L82").Value = 9068.429
$ws.Range("M82").Value = -1096.4546
$ws.Range("N82").Value = -9790.429
$ws.Range("H85").Value = 4417.278
$ws.Range("I85").Value = 1457.4546
$ws.Range("J85").Value = 9068.429
$ws.Range("K85").Value = 1457.4546
$ws.Range("L85").Value = 9068.429
$ws.Range("M85").Value = -209.4546
$ws.Range("N85").Value = -11564.429
$ws.Range("H122").Value = 18853818
$ws.Range("I122").Value = 13566323
$ws.Range("J122").Value = 40003800
$ws.Range("K122").Value = 40698969
$ws.Range("L122").Value = 120011400
$ws.Range("M122").Value = -40696519
$ws.Range("N122").Value = -120016300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 41572
$ws.Range("J138").Value = 41572
$ws.Range("L138").Value = 41572
$ws.Range("N138").Value = -51852
